$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhanCongNganHan")

# Row 5 (task 3.3 Quan Ly Trong Tai - Nhat Anh): add interim + final results
# (written first so the new shared string "100% (9/06/2010)" gets inserted
# before "100% (8/06/2010)", matching the target shared-string order)
$ws.Range("D5").Value = "70% (Xong 2/3 Chức năng)"
$ws.Range("E5").Value = "100% (9/06/2010)"

# Row 4 (task 3.2 - Tuan Anh): add final result in column E
$ws.Range("E4").Value = "100% (9/06/2010)"

# Row 3 (task 3.1 - Van Hoang): add final result in column E
$ws.Range("E3").Value = "100% (8/06/2010)"

# Update view: scroll back to A1 and move selection to E6
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E6").Select()
